# "fixed names on allgraph"
# Add the first two racers' rows (MoscaMye / Zokalyx) to every sheet of the
# scoreboard, matching the style already used for the header-adjacent data
# (centered horizontal alignment == the workbook's existing style index 1),
# and leave a gray-filled placeholder cell at F3 (Race 2 result still TBD
# for Zokalyx).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Reserve the centered-alignment style (same as the rest of the data
    # rows) across A2:F2 and A3:E3 before writing values. F3 is left out on
    # purpose - it only gets a fill, not the centered alignment.
    $ws.Range("A2:F2").HorizontalAlignment = -4108
    $ws.Range("A3:E3").HorizontalAlignment = -4108

    # Row 2 - MoscaMye, 1st place, 200 points
    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = "MoscaMye"
    $ws.Range("C2").Value = "1st"
    $ws.Range("D2").Value = 200
    $ws.Range("E2").Value = "1st (+100)"
    $ws.Range("F2").Value = "1st (+100)"

    # Row 3 - Zokalyx, 2nd place, 99 points
    $ws.Range("A3").Value = 27
    $ws.Range("B3").Value = "Zokalyx"
    $ws.Range("C3").Value = "2nd"
    $ws.Range("D3").Value = 99
    $ws.Range("E3").Value = "2nd (+99)"

    # F3 stays empty (Race 2 not yet run for Zokalyx) but gets a light gray
    # fill to mark it.
    $ws.Range("F3").Interior.Color = 13421772
}
